# Fix gitignore in exercise
# Applies the changes described by the commit:
#  - bump the "last updated" date on the slide master from 06.06.2024 to 10.06.2024
#  - use "**/" (double-star) glob patterns instead of "*/" in the .gitignore example (slide 6)
#  - clarify wording on slide 8 ("... hier leer bzw. nicht vorhanden sein.")
#  - clarify wording + extra .gitignore file on slide 9

$p = $ppt.ActivePresentation

function Replace-Substring($TextRange, $Old, $New) {
    $full = $TextRange.Text
    $idx = $full.IndexOf($Old)
    if ($idx -lt 0) {
        throw "Substring not found: [$Old]"
    }
    $chars = $TextRange.Characters($idx + 1, $Old.Length)
    $chars.Text = $New
    return $chars
}

# ---------------------------------------------------------------------------
# 1) Slide master date placeholder: 06.06.2024 -> 10.06.2024
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t -eq "06.06.2024") {
            Replace-Substring $shp.TextFrame.TextRange "06.06.2024" "10.06.2024" | Out-Null
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 6 (.gitignore exercise): "*/" -> "**/" patterns
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$body6 = $slide6.Shapes.Item(1).TextFrame.TextRange

# First line of the .gitignore example is its own run "*/" immediately
# followed by a separate "build" run and a separate "/" run - only touch
# the "*/" run itself so "build" (err="1") and "/" stay intact.
$full6 = $body6.Text
$idxBuild = $full6.IndexOf("*/build/")
$body6.Characters($idxBuild + 1, 2).Text = "**/"

# Second line "*/bin/" is a single run - replace it wholesale.
Replace-Substring $body6 "*/bin/" "**/bin/" | Out-Null

# ---------------------------------------------------------------------------
# 3) Slide 8 (Branches exercise): clarify wording
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$body8 = $slide8.Shapes.Item(1).TextFrame.TextRange

Replace-Substring $body8 " hier leer sein." " hier leer bzw. nicht vorhanden sein." | Out-Null

# ---------------------------------------------------------------------------
# 4) Slide 9 (Branches exercise continued)
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$body9 = $slide9.Shapes.Item(1).TextFrame.TextRange

# 4a) "Legen Sie im Ordner " -> "Legen Sie erneut " + "den Ordner " (two runs)
Replace-Substring $body9 "Legen Sie im Ordner " "Legen Sie erneut den Ordner " | Out-Null
# re-split the just-written text into its two runs
$full9 = $body9.Text
$splitAt = $full9.IndexOf("Legen Sie erneut den Ordner ")
$part1 = "Legen Sie erneut "
$body9.Characters($splitAt + 1, $part1.Length).Text = $part1

# 4b) " eine Datei " -> ", " (red) + "sowie eine Datei "
Replace-Substring $body9 " eine Datei " ", sowie eine Datei " | Out-Null
$full9 = $body9.Text
$commaAt = $full9.IndexOf(", sowie eine Datei ")
$commaRange = $body9.Characters($commaAt + 1, 2)
$commaRange.Text = ", "
$commaRange.Font.Color.RGB = 192   # RRGGBB C00000 -> COM RGB 0x0000C0 == 192

# 4c) " Ordner ist hier immer noch leer." -> " Ordner ist erneut nicht vorhanden oder leer."
Replace-Substring $body9 " Ordner ist hier immer noch leer." " Ordner ist erneut nicht vorhanden oder leer." | Out-Null
